$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the old "_GoBack" bookmark used to sit between "ca" and
# "tegories" in the first paragraph, splitting one sentence into two runs.
# The edited document no longer has a split there - it's a single run.
# We delete that span and retype it so it collapses back into one run
# (the bookmark that lived inside the deleted span disappears with it),
# then nudge formatting on/off to stop the engine from re-merging the new
# run with the following (unrelated, pre-existing) run that holds just a
# trailing space.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Without further analysis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start

$r2 = $d.Content
$r2.Find.Execute("launch dates.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$end = $r2.End

$full = $d.Range($start, $end)
$savedText = $full.Text
$full.Delete()

$pt = $d.Range($start, $start)
$pt.InsertAfter($savedText)

$newRange = $d.Range($start, $start + $savedText.Length)
$newRange.Font.Bold = 1
$newRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: "... canceled or could provide ..." -> "... canceled could
# provide ..." (drop the stray " or").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("canceled or could provide", $true, $false, $false, $false, $false, `
    $true, 1, $false, "canceled could provide", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: "... whereas a significant number of ..." -> "... whereas a
# significant proportion of ...", with "proportion" ending up as its own
# run (the surrounding text stays in separate runs too, matching the
# authored edit where the word was retyped independently).
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("whereas a significant number of ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$sentenceStart = $r3.Start
$sentenceEnd = $r3.End

$numRange = $d.Range($sentenceStart, $sentenceEnd)
$numRange.Find.Execute("number", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$numStart = $numRange.Start
$numEnd = $numRange.End

$wordRange = $d.Range($numStart, $numEnd)
$wordRange.Text = "proportion"
$newWordEnd = $numStart + "proportion".Length
$newSentenceEnd = $sentenceEnd - ($numEnd - $numStart) + "proportion".Length

# Isolate "proportion" into its own run.
$isolate = $d.Range($numStart, $newWordEnd)
$isolate.Font.Bold = 1
$isolate.Font.Bold = 0

# Isolate the trailing " of " into its own run too, so it doesn't stay
# merged with "proportion" (the following "plays" is already its own
# italic run, so this range naturally ends right before it).
$afterRange = $d.Range($newWordEnd, $newSentenceEnd)
$afterRange.Font.Bold = 1
$afterRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Change 4: the "_GoBack" bookmark reappears later in the document, now
# splitting "between" into "bet" / "ween" inside "... a correlation
# between success rate ...".
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("correlation between", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$corrStart = $r4.Start
$splitPos = $corrStart + "correlation bet".Length

$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null
